$wb = $excel.ActiveWorkbook

# Update "想去人数" (interested-count) figures on both the "展览" and
# "全部类型" sheets, which carry duplicate data tables.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 547
    $ws.Range("F4").Value = 272
    $ws.Range("F7").Value = 783
}
